$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update F column "想去人数" (want-to-go count) values
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 197
$wsExpo.Range("F3").Value = 525
$wsExpo.Range("F4").Value = 37
$wsExpo.Range("F7").Value = 34
$wsExpo.Range("F8").Value = 27
$wsExpo.Range("F9").Value = 267
$wsExpo.Range("F10").Value = 2870

# Sheet "全部类型" (All types) - update F column "想去人数" (want-to-go count) values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 197
$wsAll.Range("F4").Value = 525
$wsAll.Range("F5").Value = 37
$wsAll.Range("F8").Value = 34
$wsAll.Range("F9").Value = 27
$wsAll.Range("F10").Value = 267
$wsAll.Range("F11").Value = 2870
